$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Fix the "Generierung neuer Services..." bullet: merge the two runs
#    that were split by the stray "_GoBack" bookmark into one clean run
#    with the full sentence, and drop that bookmark (it gets re-created
#    further down, in the newly added body text).
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$r7 = $p7.Range
$r7.MoveEnd(1, -1)
$r7.Text = "Generierung neuer Services wird durch die Erstellung neuer Technologien möglich"

# ---------------------------------------------------------------------
# 2. Helper: fill an (empty) paragraphs content by inserting a raw
#    WordprocessingML <w:p> fragment into its range (mark excluded), so
#    runs/bookmarks/proofErr markers come out exactly as specified,
#    without Word silently merging adjacent runs.
# ---------------------------------------------------------------------
function Fill-Paragraph($idx, $bodyXml) {
    $p = $d.Paragraphs($idx)
    $r = $p.Range
    $r.MoveEnd(1, -1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 3. Append 6 new empty paragraphs after the last paragraph currently in
#    the document (the blank one that follows the "Zitat" quote). Insert
#    one extra scratch paragraph too, so none of the 6 we care about is
#    ever the documents very last paragraph mark while we fill it in
#    (that special position behaves differently under InsertXML).
# ---------------------------------------------------------------------
$pAnchor = $d.Paragraphs(10)
$rngAnchor = $pAnchor.Range
for ($i = 0; $i -lt 7; $i++) {
    $rngAnchor.InsertParagraphAfter()
}

Fill-Paragraph 11 '<w:p><w:r><w:t>Der Begriff Internet der Dinge beschreibt in unseren Augen, dass die Technik zunehmend in den Hintergrund rücken wird, bis sie schlussendlich scheinbar verschwindet und herkömmlichen Gegenständen Intelligenz verleiht. Die Aufmerksamkeit der Nutzer soll entgegen dem momentanen Stand nicht auf der Bedienung der Geräte liegen</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Vielmehr sollen Endgeräte</w:t></w:r><w:r><w:t xml:space="preserve">, die Interaktionsmöglichkeiten  </w:t></w:r><w:r><w:t xml:space="preserve">bieten, den Endanwender bei seinen Tätigkeiten unterstützen und nicht dessen volle Aufmerksamkeit einfordern. </w:t></w:r></w:p>'
Fill-Paragraph 12 '<w:p><w:r><w:t xml:space="preserve">Ein Ziel der Thematik ist es, die Lücke zwischen realer und virtueller Welt zu minimieren. Hierzu werden eindeutig identifizierbare reale Objekte virtuell abgebildet und erhalten die Möglichkeit Informationen entweder mit uns oder auch primär mit anderen Geräten teilen zu können. </w:t></w:r><w:r><w:t>Durch die stetige Weiterentwickl</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ung der Technologien wird sich die Anzahl an Informationsquellen, die auf Bedarf abgerufen werden können, rapide ansteigend</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">vergrößern, wodurch </w:t></w:r><w:r><w:t xml:space="preserve">immense Möglichkeiten </w:t></w:r><w:r><w:t>für neue Services auftun.</w:t></w:r></w:p>'
Fill-Paragraph 13 '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Internet </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>of</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> Food</w:t></w:r></w:p>'
Fill-Paragraph 14 '<w:p><w:r><w:t xml:space="preserve">Per Definition ist das Internet der Lebensmittel größtenteils auf Produktionskette zugeschnitten und nicht auf den Konsumenten. Informationen über Haltbarkeit, Nährwerte, Inhaltsstoffe und Herkunft gelangen im Regelfall nur in gedruckter Form auf den Verpackungen zum Käufer. Ergänzt werden diese Informationen eventuell noch über Webseiten der Hersteller. </w:t></w:r></w:p>'
Fill-Paragraph 15 '<w:p><w:r><w:t xml:space="preserve">An diesen Grenzen Knüpfen wir an und zeigen einen Umsetzung, die zeitnah realisierbar ist, für die Übertragung und Verknüpfung der Lebensmittel mit dem Internet auf. Wir lassen die Lebensmittel ein Teil des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IoT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> werden und kreieren die Möglichkeit Lebensmittel kommunizieren zu lassen. Eine Brücke zwischen Einkaufsmärkten, dem digitalen intelligenten Kühlschrank und einem Netzwerk, das sich mit dem nachhaltigen Umgang mit Lebensmitteln auseinandersetzt, wird geschaffen.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
Fill-Paragraph 16 '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>'

# ---------------------------------------------------------------------
# 4. Drop the 17th (scratch) paragraph again by deleting its paragraph
#    mark; this merges it away without disturbing paragraph 16s own
#    (bold) formatting, leaving 16 as the documents final paragraph.
# ---------------------------------------------------------------------
$p16 = $d.Paragraphs(16)
$p17 = $d.Paragraphs(17)
$delRng = $d.Range($p16.Range.End, $p17.Range.End)
$delRng.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
